$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (UsqWnbLa): odds refresh ---
$ws.Range("O2").Value = 1.5
$ws.Range("P2").Value = 2.5
$ws.Range("Q2").Value = 2.6
$ws.Range("R2").Value = 1.48

# --- Row 7 (vkSujNA4): odds refresh ---
$ws.Range("G7").Value = 2.05
$ws.Range("I7").Value = 3.3
$ws.Range("J7").Value = 2.6
$ws.Range("K7").Value = 2.38
$ws.Range("U7").Value = 1.53
$ws.Range("V7").Value = 2.38
$ws.Range("X7").Value = 12
$ws.Range("AD7").Value = 7.5
$ws.Range("AF7").Value = 34
$ws.Range("AK7").Value = 34
$ws.Range("AM7").Value = 26
$ws.Range("AP7").Value = 17
$ws.Range("AR7").Value = 41
$ws.Range("AU7").Value = 7
$ws.Range("AY7").Value = 21

# --- Row 8 (hGAStJXj): odds refresh ---
$ws.Range("M8").Value = 1.04
$ws.Range("O8").Value = 1.2

# --- Row 11 (QLdg8mmJ): kickoff time moved ---
$ws.Range("C11").Value = "07:30"

# --- Row 12 (C8BSPIJt): odds refresh ---
$ws.Range("Q12").Value = 2
$ws.Range("R12").Value = 1.85

# --- Row 14 (pANvoeuf / SINGAPORE - Young Lions vs Albirex Niigata) removed entirely;
# remaining fixtures shift up one row ---
$ws.Rows(14).Delete()

# --- (was row 17, now row 16 - EkuvK9MU): odds refresh ---
$ws.Range("G16").Value = 2.35
$ws.Range("H16").Value = 3.25
$ws.Range("I16").Value = 3.1
$ws.Range("J16").Value = 3
$ws.Range("K16").Value = 2.1
$ws.Range("M16").Value = 1.06
$ws.Range("N16").Value = 10
$ws.Range("O16").Value = 1.29
$ws.Range("P16").Value = 3.5
$ws.Range("Q16").Value = 2
$ws.Range("R16").Value = 1.85
$ws.Range("S16").Value = 1.4
$ws.Range("T16").Value = 2.75
$ws.Range("U16").Value = 1.73
$ws.Range("V16").Value = 2
$ws.Range("W16").Value = 8.5
$ws.Range("AA16").Value = 19
$ws.Range("AC16").Value = 10
$ws.Range("AD16").Value = 6
$ws.Range("AG16").Value = 201
$ws.Range("AH16").Value = 10
$ws.Range("AI16").Value = 15
$ws.Range("AJ16").Value = 11
$ws.Range("AM16").Value = 34
$ws.Range("AO16").Value = 13
$ws.Range("AP16").Value = 23
$ws.Range("AR16").Value = 67
$ws.Range("AT16").Value = 2.75
$ws.Range("AU16").Value = 8
$ws.Range("AY16").Value = 26

# --- (was row 19, now row 18 - G6D3bU7d): odds refresh ---
$ws.Range("O18").Value = 1.14
$ws.Range("P18").Value = 5.86
$ws.Range("S18").Value = 1.26
$ws.Range("T18").Value = 3.92
$ws.Range("U18").Value = 2.98
$ws.Range("V18").Value = 1.32

# --- (was row 20, now row 19 - OjaSazmD): odds refresh ---
$ws.Range("G19").Value = 2.55
$ws.Range("H19").Value = 2.82
$ws.Range("I19").Value = 2.92
$ws.Range("J19").Value = 3.2
$ws.Range("L19").Value = 3.45
$ws.Range("M19").Value = 1.08
$ws.Range("N19").Value = 7.75
$ws.Range("V19").Value = 1.8
$ws.Range("W19").Value = 6.9
$ws.Range("X19").Value = 11.75
$ws.Range("Z19").Value = 29
$ws.Range("AB19").Value = 37
$ws.Range("AD19").Value = 5.5
$ws.Range("AF19").Value = 75
$ws.Range("AH19").Value = 8
$ws.Range("AI19").Value = 14.5
$ws.Range("AL19").Value = 27
$ws.Range("AM19").Value = 37
$ws.Range("AP19").Value = 23
$ws.Range("AR19").Value = 110
$ws.Range("AU19").Value = 6.8
$ws.Range("AW19").Value = 4.7
$ws.Range("AY19").Value = 23
$ws.Range("AZ19").Value = 75
$ws.Range("BA19").Value = 110
$ws.Range("BB19").Value = 300
